$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '29.620.77'
$ws.Cells.Item(2, 5).Value2 = '  -2.23%  '
Set-TextValue 3 4 '2.003.55'
$ws.Cells.Item(3, 5).Value2 = '  -4.43%  '
$ws.Cells.Item(4, 5).Value2 = '  +0.45%  '
Set-TextValue 5 4 '330.61'
$ws.Cells.Item(5, 5).Value2 = '  -3.92%  '
Set-TextValue 6 4 '1.012'
$ws.Cells.Item(6, 5).Value2 = '  +0.55%  '
Set-TextValue 7 4 '0.5017'
$ws.Cells.Item(7, 5).Value2 = '  -3.77%  '
Set-TextValue 8 4 '0.4245'
$ws.Cells.Item(8, 5).Value2 = '  -3.59%  '
Set-TextValue 9 4 '54.26'
$ws.Cells.Item(9, 5).Value2 = '  -0.27%  '
Set-TextValue 10 4 '0.08963'
$ws.Cells.Item(10, 5).Value2 = '  -3.32%  '
Set-TextValue 11 4 '1.122'
$ws.Cells.Item(11, 5).Value2 = '  -3.93%  '
Set-TextValue 12 4 '23.44'
$ws.Cells.Item(12, 5).Value2 = '  -5.01%  '
Set-TextValue 13 4 '2.037.33'
$ws.Cells.Item(13, 5).Value2 = '  -2.42%  '
Set-TextValue 14 4 '8.096'
$ws.Cells.Item(14, 5).Value2 = '  -6.25%  '
Set-TextValue 15 4 '6.504'
$ws.Cells.Item(15, 5).Value2 = '  -5.59%  '
$ws.Cells.Item(16, 5).Value2 = '  +0.49%  '
Set-TextValue 17 4 '94.58'
$ws.Cells.Item(17, 5).Value2 = '  -6.80%  '
$ws.Cells.Item(18, 5).Value2 = '  -3.56%  '
Set-TextValue 19 4 '0.06679'
$ws.Cells.Item(19, 5).Value2 = '  -0.44%  '
Set-TextValue 20 4 '19.80'
$ws.Cells.Item(20, 5).Value2 = '  -6.15%  '
Set-TextValue 21 4 '1.012'
$ws.Cells.Item(21, 5).Value2 = '  +0.60%  '
Set-TextValue 22 4 '5.961'
$ws.Cells.Item(22, 5).Value2 = '  -6.04%  '
Set-TextValue 23 4 '29.607.14'
$ws.Cells.Item(23, 5).Value2 = '  -2.53%  '
Set-TextValue 24 4 '12.00'
$ws.Cells.Item(24, 5).Value2 = '  -3.96%  '
$ws.Cells.Item(25, 5).Value2 = '  -0.72%  '
Set-TextValue 26 4 '159.64'
$ws.Cells.Item(26, 5).Value2 = '  -1.71%  '
Set-TextValue 27 4 '20.74'
$ws.Cells.Item(27, 5).Value2 = '  -4.97%  '
Set-TextValue 28 4 '6.323'
$ws.Cells.Item(28, 5).Value2 = '  -5.58%  '
Set-TextValue 29 4 '2.307'
$ws.Cells.Item(29, 5).Value2 = '  -7.83%  '
Set-TextValue 30 4 '128.56'
$ws.Cells.Item(30, 5).Value2 = '  -3.42%  '
Set-TextValue 31 4 '1.059'
$ws.Cells.Item(31, 5).Value2 = '  -6.31%  '
Set-TextValue 32 4 '0.09957'
$ws.Cells.Item(32, 5).Value2 = '  -5.07%  '
Set-TextValue 33 4 '1.565'
$ws.Cells.Item(33, 5).Value2 = '  -5.46%  '
Set-TextValue 34 4 '5.858'
$ws.Cells.Item(34, 5).Value2 = '  -5.63%  '
Set-TextValue 35 4 '3.808'
$ws.Cells.Item(35, 5).Value2 = '  -2.88%  '
Set-TextValue 36 4 '9.456'
$ws.Cells.Item(36, 5).Value2 = '  -7.26%  '
Set-TextValue 37 4 '0.02473'
$ws.Cells.Item(37, 5).Value2 = '  -5.74%  '
Set-TextValue 38 4 '1.315'
$ws.Cells.Item(38, 5).Value2 = '  -2.17%  '
Set-TextValue 39 4 '0.06369'
$ws.Cells.Item(39, 5).Value2 = '  -5.68%  '
Set-TextValue 40 4 '0.6578'
$ws.Cells.Item(40, 5).Value2 = '  -5.71%  '
Set-TextValue 41 4 '11.73'
$ws.Cells.Item(41, 5).Value2 = '  -6.06%  '
Set-TextValue 42 4 '0.2062'
$ws.Cells.Item(42, 5).Value2 = '  -6.81%  '
$ws.Cells.Item(43, 5).Value2 = '  +0.54%  '
Set-TextValue 44 4 '0.6348'
$ws.Cells.Item(44, 5).Value2 = '  -6.55%  '
Set-TextValue 45 4 '13.45'
$ws.Cells.Item(45, 5).Value2 = '  -6.19%  '
Set-TextValue 46 4 '2.212'
$ws.Cells.Item(46, 5).Value2 = '  -5.36%  '
Set-TextValue 47 4 '1.317'
$ws.Cells.Item(47, 5).Value2 = '  -4.75%  '
Set-TextValue 48 4 '3.523'
$ws.Cells.Item(48, 5).Value2 = '  -3.34%  '
Set-TextValue 49 4 '0.00000000338'
$ws.Cells.Item(49, 5).Value2 = '  -3.71%  '
Set-TextValue 50 4 '0.06993'
$ws.Cells.Item(50, 5).Value2 = '  -3.33%  '
Set-TextValue 51 4 '1.128'
$ws.Cells.Item(51, 5).Value2 = '  -6.59%  '
